$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.359913
$ws.Range("H2").Value = 31.079739
$ws.Range("I2").Value = 0.2499874361758538
$ws.Range("J2").Value = 0.2499874361758538
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.44605333333334
$ws.Range("N2").Value = 94.33816
$ws.Range("O2").Value = 0.273208187120734
$ws.Range("P2").Value = 0.273208187120734
$ws.Range("Q2").Value = 325.7783767266933
$ws.Range("R2").Value = 2932.00539054024
$ws.Range("S2").Value = 0.06829861424056521
$ws.Range("T2").Value = 0.06829861424056521

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.359913
$ws.Range("H3").Value = 31.079739
$ws.Range("I3").Value = 0.2499874361758538
$ws.Range("J3").Value = 0.2499874361758538
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.80064033333333
$ws.Range("N3").Value = 173.401921
$ws.Range("O3").Value = 0.5021809252974907
$ws.Range("P3").Value = 0.5021809252974908
$ws.Range("Q3").Value = 598.8096051976243
$ws.Range("R3").Value = 5389.286446778618
$ws.Range("S3").Value = 0.1255389220115377
$ws.Range("T3").Value = 0.1255389220115377

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.359913
$ws.Range("H4").Value = 31.079739
$ws.Range("I4").Value = 0.2499874361758538
$ws.Range("J4").Value = 0.2499874361758538
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.04190666666667
$ws.Range("N4").Value = 36.12572
$ws.Range("O4").Value = 0.1046219522368387
$ws.Range("P4").Value = 0.1046219522368387
$ws.Range("Q4").Value = 124.7531054207867
$ws.Range("R4").Value = 1122.77794878708
$ws.Range("S4").Value = 0.02615417360739993
$ws.Range("T4").Value = 0.02615417360739993

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.359913
$ws.Range("H5").Value = 31.079739
$ws.Range("I5").Value = 0.2499874361758538
$ws.Range("J5").Value = 0.2499874361758538
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.81063466666667
$ws.Range("N5").Value = 41.431904
$ws.Range("O5").Value = 0.1199889353449366
$ws.Range("P5").Value = 0.1199889353449366
$ws.Range("Q5").Value = 143.0769736214506
$ws.Range("R5").Value = 1287.692762593056
$ws.Range("S5").Value = 0.02999572631635099
$ws.Range("T5").Value = 0.02999572631635099

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.857753
$ws.Range("H6").Value = 38.57325899999999
$ws.Range("I6").Value = 0.3102609749186496
$ws.Range("J6").Value = 0.3102609749186496
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.44605333333334
$ws.Range("N6").Value = 94.33816
$ws.Range("O6").Value = 0.273208187120734
$ws.Range("P6").Value = 0.273208187120734
$ws.Range("Q6").Value = 404.3255865848266
$ws.Range("R6").Value = 3638.930279263439
$ws.Range("S6").Value = 0.08476583849183579
$ws.Range("T6").Value = 0.08476583849183579

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.857753
$ws.Range("H7").Value = 38.57325899999999
$ws.Range("I7").Value = 0.3102609749186496
$ws.Range("J7").Value = 0.3102609749186496
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.80064033333333
$ws.Range("N7").Value = 173.401921
$ws.Range("O7").Value = 0.5021809252974907
$ws.Range("P7").Value = 0.5021809252974908
$ws.Range("Q7").Value = 743.1863566478374
$ws.Range("R7").Value = 6688.677209830537
$ws.Range("S7").Value = 0.155807143468349
$ws.Range("T7").Value = 0.1558071434683491

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.857753
$ws.Range("H8").Value = 38.57325899999999
$ws.Range("I8").Value = 0.3102609749186496
$ws.Range("J8").Value = 0.3102609749186496
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.04190666666667
$ws.Range("N8").Value = 36.12572
$ws.Range("O8").Value = 0.1046219522368387
$ws.Range("P8").Value = 0.1046219522368387
$ws.Range("Q8").Value = 154.8318615690533
$ws.Range("R8").Value = 1393.48675412148
$ws.Range("S8").Value = 0.03246010889889396
$ws.Range("T8").Value = 0.03246010889889396

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.857753
$ws.Range("H9").Value = 38.57325899999999
$ws.Range("I9").Value = 0.3102609749186496
$ws.Range("J9").Value = 0.3102609749186496
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.81063466666667
$ws.Range("N9").Value = 41.431904
$ws.Range("O9").Value = 0.1199889353449366
$ws.Range("P9").Value = 0.1199889353449366
$ws.Range("Q9").Value = 177.5737293172373
$ws.Range("R9").Value = 1598.163563855136
$ws.Range("S9").Value = 0.03722788405957085
$ws.Range("T9").Value = 0.03722788405957085

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 13.70542966666667
$ws.Range("H10").Value = 41.116289
$ws.Range("I10").Value = 0.3307156367103167
$ws.Range("J10").Value = 0.3307156367103167
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.44605333333334
$ws.Range("N10").Value = 94.33816
$ws.Range("O10").Value = 0.273208187120734
$ws.Range("P10").Value = 0.273208187120734
$ws.Range("Q10").Value = 430.9816722542489
$ws.Range("R10").Value = 3878.83505028824
$ws.Range("S10").Value = 0.0903542195581049
$ws.Range("T10").Value = 0.0903542195581049

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 13.70542966666667
$ws.Range("H11").Value = 41.116289
$ws.Range("I11").Value = 0.3307156367103167
$ws.Range("J11").Value = 0.3307156367103167
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 57.80064033333333
$ws.Range("N11").Value = 173.401921
$ws.Range("O11").Value = 0.5021809252974907
$ws.Range("P11").Value = 0.5021809252974908
$ws.Range("Q11").Value = 792.1826107767965
$ws.Range("R11").Value = 7129.643496991169
$ws.Range("S11").Value = 0.1660790844535356
$ws.Range("T11").Value = 0.1660790844535356

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 13.70542966666667
$ws.Range("H12").Value = 41.116289
$ws.Range("I12").Value = 0.3307156367103167
$ws.Range("J12").Value = 0.3307156367103167
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.04190666666667
$ws.Range("N12").Value = 36.12572
$ws.Range("O12").Value = 0.1046219522368387
$ws.Range("P12").Value = 0.1046219522368387
$ws.Range("Q12").Value = 165.0395048725645
$ws.Range("R12").Value = 1485.35554385308
$ws.Range("S12").Value = 0.03460011554788244
$ws.Range("T12").Value = 0.03460011554788244

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 13.70542966666667
$ws.Range("H13").Value = 41.116289
$ws.Range("I13").Value = 0.3307156367103167
$ws.Range("J13").Value = 0.3307156367103167
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.81063466666667
$ws.Range("N13").Value = 41.431904
$ws.Range("O13").Value = 0.1199889353449366
$ws.Range("P13").Value = 0.1199889353449366
$ws.Range("Q13").Value = 189.2806820760284
$ws.Range("R13").Value = 1703.526138684256
$ws.Range("S13").Value = 0.03968221715079374
$ws.Range("T13").Value = 0.03968221715079374

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.518638999999999
$ws.Range("H14").Value = 13.555917
$ws.Range("I14").Value = 0.10903595219518
$ws.Range("J14").Value = 0.10903595219518
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.44605333333334
$ws.Range("N14").Value = 94.33816
$ws.Range("O14").Value = 0.273208187120734
$ws.Range("P14").Value = 0.273208187120734
$ws.Range("Q14").Value = 142.09336298808
$ws.Range("R14").Value = 1278.84026689272
$ws.Range("S14").Value = 0.02978951483022814
$ws.Range("T14").Value = 0.02978951483022814

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.518638999999999
$ws.Range("H15").Value = 13.555917
$ws.Range("I15").Value = 0.10903595219518
$ws.Range("J15").Value = 0.10903595219518
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 57.80064033333333
$ws.Range("N15").Value = 173.401921
$ws.Range("O15").Value = 0.5021809252974907
$ws.Range("P15").Value = 0.5021809252974908
$ws.Range("Q15").Value = 261.1802276351729
$ws.Range("R15").Value = 2350.622048716556
$ws.Range("S15").Value = 0.05475577536406845
$ws.Range("T15").Value = 0.05475577536406846

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.518638999999999
$ws.Range("H16").Value = 13.555917
$ws.Range("I16").Value = 0.10903595219518
$ws.Range("J16").Value = 0.10903595219518
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.04190666666667
$ws.Range("N16").Value = 36.12572
$ws.Range("O16").Value = 0.1046219522368387
$ws.Range("P16").Value = 0.1046219522368387
$ws.Range("Q16").Value = 54.41302909836
$ws.Range("R16").Value = 489.71726188524
$ws.Range("S16").Value = 0.01140755418266234
$ws.Range("T16").Value = 0.01140755418266234

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.518638999999999
$ws.Range("H17").Value = 13.555917
$ws.Range("I17").Value = 0.10903595219518
$ws.Range("J17").Value = 0.10903595219518
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.81063466666667
$ws.Range("N17").Value = 41.431904
$ws.Range("O17").Value = 0.1199889353449366
$ws.Range("P17").Value = 0.1199889353449366
$ws.Range("Q17").Value = 177.5737293172373
$ws.Range("R17").Value = 1598.163563855136
$ws.Range("S17").Value = 0.03722788405957085
$ws.Range("T17").Value = 0.03722788405957085
